$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.067.81"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "1.928.58"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.86"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4608"
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3831"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07753"
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9807"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.58"
$ws.Range("E11").Value = "  +3.08%  "
$ws.Range("D12").Value = "1.947.60"
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.982"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.695"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07030"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.45"
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009568"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.78"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "29.107.27"
$ws.Range("E21").Value = "  +1.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.351"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.97"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.078"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.71"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.08"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.680"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.88"
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.848"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09325"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8656"
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.134"
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.250"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.159"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02050"
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.057"
$ws.Range("E39").Value = "  +13.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.549"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5527"
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1754"
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.376"
$ws.Range("E43").Value = "  +2.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000002878"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.215"
$ws.Range("E45").Value = "  +6.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5211"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.24"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06932"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.782"
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.49"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").Value = "  +0.15%  "
